$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.020519786208548
$ws.Range("D2").Value = 1.032980814537891
$ws.Range("E2").Value = 1.030862505864558
$ws.Range("F2").Value = 1.040391049643465
$ws.Range("I2").Value = 1.032242934973097
$ws.Range("J2").Value = 1.025716435313205
$ws.Range("K2").Value = 1.035784480774071
$ws.Range("L2").Value = 1.033672288707147
$ws.Range("M2").Value = 1.043173531642934
$ws.Range("N2").Value = 1.012598234192461
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.021505472363392
$ws.Range("D3").Value = 1.033373492209532
$ws.Range("E3").Value = 1.031742273614255
$ws.Range("F3").Value = 1.041367336713423
$ws.Range("I3").Value = 1.032284454184902
$ws.Range("J3").Value = 1.026339155371382
$ws.Range("K3").Value = 1.035987156426971
$ws.Range("L3").Value = 1.034360308545783
$ws.Range("M3").Value = 1.043959796009926
$ws.Range("N3").Value = 1.012809391485004
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.022143964655388
$ws.Range("D4").Value = 1.033627238479311
$ws.Range("E4").Value = 1.032312414024458
$ws.Range("F4").Value = 1.041999742801965
$ws.Range("I4").Value = 1.032309714103485
$ws.Range("J4").Value = 1.026742215390544
$ws.Range("K4").Value = 1.036117218494523
$ws.Range("L4").Value = 1.034805740158838
$ws.Range("M4").Value = 1.044468628855964
$ws.Range("N4").Value = 1.012945932868433
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.02241255076789
$ws.Range("D5").Value = 1.033733829183337
$ws.Range("E5").Value = 1.0325523086064
$ws.Range("F5").Value = 1.042265768498901
$ws.Range("I5").Value = 1.032319948025109
$ws.Range("J5").Value = 1.026911689447251
$ws.Range("K5").Value = 1.036171635829165
$ws.Range("L5").Value = 1.034993055356412
$ws.Range("M5").Value = 1.044682557038304
$ws.Range("N5").Value = 1.013003312605867
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.022457657154381
$ws.Range("D6").Value = 1.033751721226819
$ws.Range("E6").Value = 1.032592600099396
$ws.Range("F6").Value = 1.042310444854889
$ws.Range("I6").Value = 1.032321643721763
$ws.Range("J6").Value = 1.026940146485414
$ws.Range("K6").Value = 1.036180757402142
$ws.Range("L6").Value = 1.035024509644194
$ws.Range("M6").Value = 1.044718477356461
$ws.Range("N6").Value = 1.013012945600639
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.022147552872711
$ws.Range("D7").Value = 1.033628663083328
$ws.Range("E7").Value = 1.032315618693528
$ws.Range("F7").Value = 1.042003296814807
$ws.Range("I7").Value = 1.032309852364964
$ws.Range("J7").Value = 1.026744479803571
$ws.Range("K7").Value = 1.036117946648213
$ws.Range("L7").Value = 1.034808242855803
$ws.Range("M7").Value = 1.044471487318009
$ws.Range("N7").Value = 1.012946699667204
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.020852760745264
$ws.Range("D8").Value = 1.033113591353271
$ws.Range("E8").Value = 1.031159646269525
$ws.Range("F8").Value = 1.040720848181848
$ws.Range("I8").Value = 1.032257298522637
$ws.Range("J8").Value = 1.02592686148736
$ws.Range("K8").Value = 1.03585319918219
$ws.Range("L8").Value = 1.033904758644544
$ws.Range("M8").Value = 1.043439238701704
$ws.Range("N8").Value = 1.01266961452193
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.018576463991058
$ws.Range("D9").Value = 1.032203445348818
$ws.Range("E9").Value = 1.029129404054509
$ws.Range("F9").Value = 1.038466300984338
$ws.Range("I9").Value = 1.032152433975761
$ws.Range("J9").Value = 1.024487053994128
$ws.Range("K9").Value = 1.035378457918984
$ws.Range("L9").Value = 1.032314558708309
$ws.Range("M9").Value = 1.041620854305389
$ws.Range("N9").Value = 1.012180669879588
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.017062526554289
$ws.Range("D10").Value = 1.031595118061086
$ws.Range("E10").Value = 1.027780502003229
$ws.Range("F10").Value = 1.036966906323638
$ws.Range("I10").Value = 1.032074337206308
$ws.Range("J10").Value = 1.023527858083324
$ws.Range("K10").Value = 1.0350565289377
$ws.Range("L10").Value = 1.031255727238729
$ws.Range("M10").Value = 1.040409053809095
$ws.Range("N10").Value = 1.01185426684906
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.016407834258291
$ws.Range("D11").Value = 1.031331363122159
$ws.Range("E11").Value = 1.02719751604596
$ws.Range("F11").Value = 1.036318530867864
$ws.Range("I11").Value = 1.032038590045243
$ws.Range("J11").Value = 1.023112684820755
$ws.Range("K11").Value = 1.034915860819369
$ws.Range("L11").Value = 1.030797561586085
$ws.Range("M11").Value = 1.039884453977792
$ws.Range("N11").Value = 1.011712830674984
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.016164780977741
$ws.Range("D12").Value = 1.031233343229412
$ws.Range("E12").Value = 1.026981135108032
$ws.Range("F12").Value = 1.036077827692512
$ws.Range("I12").Value = 1.032025022792516
$ws.Range("J12").Value = 1.022958496437086
$ws.Range("K12").Value = 1.034863421092832
$ws.Range("L12").Value = 1.030627426717062
$ws.Range("M12").Value = 1.039689613111889
$ws.Range("N12").Value = 1.011660280064749
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.016216910904582
$ws.Range("D13").Value = 1.031254371035287
$ws.Range("E13").Value = 1.027027542023833
$ws.Range("F13").Value = 1.036129453331565
$ws.Range("I13").Value = 1.032027946085112
$ws.Range("J13").Value = 1.022991569206139
$ws.Range("K13").Value = 1.034874678134057
$ws.Range("L13").Value = 1.030663919026335
$ws.Range("M13").Value = 1.039731406266291
$ws.Range("N13").Value = 1.011671553018747
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.016387740772252
$ws.Range("D14").Value = 1.031323261765111
$ws.Range("E14").Value = 1.027179626537261
$ws.Range("F14").Value = 1.036298631561239
$ws.Range("I14").Value = 1.032037474465559
$ws.Range("J14").Value = 1.023099939037726
$ws.Range("K14").Value = 1.034911529988205
$ws.Range("L14").Value = 1.030783497195301
$ws.Range("M14").Value = 1.039868347979315
$ws.Range("N14").Value = 1.011708487127541
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.016493011745933
$ws.Range("D15").Value = 1.031365701118141
$ws.Range("E15").Value = 1.027273352855849
$ws.Range("F15").Value = 1.036402885426508
$ws.Range("I15").Value = 1.03204330692505
$ws.Range("J15").Value = 1.023166712660531
$ws.Range("K15").Value = 1.034934210591943
$ws.Range("L15").Value = 1.030857179672124
$ws.Range("M15").Value = 1.039952724823852
$ws.Range("N15").Value = 1.011731241484159
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.017105994331214
$ws.Range("D16").Value = 1.031612615531139
$ws.Range("E16").Value = 1.027819216079842
$ws.Range("F16").Value = 1.037009955355459
$ws.Range("I16").Value = 1.032076669026622
$ws.Range("J16").Value = 1.0235554152974
$ws.Range("K16").Value = 1.035065837980318
$ws.Range("L16").Value = 1.031286140903085
$ws.Range("M16").Value = 1.040443872353861
$ws.Range("N16").Value = 1.011863651391862
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.017490730725385
$ws.Range("D17").Value = 1.031767407624226
$ws.Range("E17").Value = 1.028161916284064
$ws.Range("F17").Value = 1.037390988733577
$ws.Range("I17").Value = 1.032097079918186
$ws.Range("J17").Value = 1.023799282817393
$ws.Range("K17").Value = 1.035148065275094
$ws.Range("L17").Value = 1.03155530199841
$ws.Range("M17").Value = 1.040751988585226
$ws.Range("N17").Value = 1.011946681686996
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.017715223185337
$ws.Range("D18").Value = 1.031857661772179
$ws.Range("E18").Value = 1.028361913336429
$ws.Range("F18").Value = 1.03761332313734
$ws.Range("I18").Value = 1.03210879889648
$ws.Range("J18").Value = 1.023941542480422
$ws.Range("K18").Value = 1.035195904350318
$ws.Range("L18").Value = 1.031712329480447
$ws.Range("M18").Value = 1.040931718855321
$ws.Range("N18").Value = 1.011995102037021
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.017791783245558
$ws.Range("D19").Value = 1.031888430358526
$ws.Range("E19").Value = 1.028430125078211
$ws.Range("F19").Value = 1.03768914764604
$ws.Range("I19").Value = 1.032112763129473
$ws.Range("J19").Value = 1.023990052019513
$ws.Range("K19").Value = 1.035212195375624
$ws.Range("L19").Value = 1.031765876918391
$ws.Range("M19").Value = 1.04099300411126
$ws.Range("N19").Value = 1.012011610449941
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.017449443642293
$ws.Range("D20").Value = 1.031750803342004
$ws.Range("E20").Value = 1.028125136841775
$ws.Range("F20").Value = 1.037350098751992
$ws.Range("I20").Value = 1.032094909291128
$ws.Range("J20").Value = 1.023773116494567
$ws.Range("K20").Value = 1.035139255746011
$ws.Range("L20").Value = 1.031526420412638
$ws.Range("M20").Value = 1.040718929453774
$ws.Range("N20").Value = 1.011937774335059
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.016337432080374
$ws.Range("D21").Value = 1.031302976512071
$ws.Range("E21").Value = 1.027134836847198
$ws.Range("F21").Value = 1.036248809125791
$ws.Range("I21").Value = 1.032034676570127
$ws.Range("J21").Value = 1.023068026115326
$ws.Range("K21").Value = 1.034900683249211
$ws.Range("L21").Value = 1.030748283065862
$ws.Range("M21").Value = 1.039828021527693
$ws.Range("N21").Value = 1.011697611362189
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.015639010498192
$ws.Range("D22").Value = 1.031021124371151
$ws.Range("E22").Value = 1.026513156935678
$ws.Range("F22").Value = 1.035557150928808
$ws.Range("I22").Value = 1.031995133291246
$ws.Range("J22").Value = 1.022624855401611
$ws.Range("K22").Value = 1.034749588660972
$ws.Range("L22").Value = 1.030259316849959
$ws.Range("M22").Value = 1.039267982108155
$ws.Range("N22").Value = 1.011546525224518
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.016009186282954
$ws.Range("D23").Value = 1.031170565835241
$ws.Range("E23").Value = 1.026842629755884
$ws.Range("F23").Value = 1.03592373896114
$ws.Range("I23").Value = 1.032016254170834
$ws.Range("J23").Value = 1.022859774364986
$ws.Range("K23").Value = 1.034829789986848
$ws.Range("L23").Value = 1.03051850029545
$ws.Range("M23").Value = 1.039564858805098
$ws.Range("N23").Value = 1.011626626899117
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.017468099243855
$ws.Range("D24").Value = 1.031758306206106
$ws.Range("E24").Value = 1.028141755561706
$ws.Range("F24").Value = 1.037368574914868
$ws.Range("I24").Value = 1.032095890679959
$ws.Range("J24").Value = 1.023784939880549
$ws.Range("K24").Value = 1.035143236772371
$ws.Range("L24").Value = 1.031539470664274
$ws.Range("M24").Value = 1.040733867417468
$ws.Range("N24").Value = 1.011941799214343
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.019164310461786
$ws.Range("D25").Value = 1.03243902517613
$ws.Range("E25").Value = 1.029653465856748
$ws.Range("F25").Value = 1.039048520278525
$ws.Range("I25").Value = 1.032180990974834
$ws.Range("J25").Value = 1.024859162826588
$ws.Range("K25").Value = 1.035502154646068
$ws.Range("L25").Value = 1.032725438074517
$ws.Range("M25").Value = 1.042090875261709
$ws.Range("N25").Value = 1.012307152852026
